$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9258052110671997
$ws.Range("B1").Value = 2.03331184387207
$ws.Range("C1").Value = 8.755983352661133
$ws.Range("D1").Value = 1.864662289619446
$ws.Range("E1").Value = 1.426128149032593
